# Generate Report for Handoff
# The build produced a new handoff package (new GUID-prefixed file names
# and new target-file hashes) and re-ran the zh-cn / de-de handoffs at
# slightly later timestamps. Update the status workbook accordingly:
#   - the source markdown file reference changed name (old GUID -> new GUID)
#   - the zh-cn / de-de handoff target file names changed (new GUID + new hash)
#   - the "Latest Handoff Datetime" for zh-cn / de-de moved forward a few seconds

$wb = $excel.ActiveWorkbook

$oldGuid = "fe7750ab-8e1a-454c-9068-0fc2e16ec1d5"
$newGuid = "a353ef61-1048-482c-864e-e5e91291b8c9"

$oldMd  = "$oldGuid.md"
$newMd  = "$newGuid.md"

$oldZh  = "$oldGuid.fddd2aa3a2d4333e1266330065820dd4a060f147.zh-cn.xlf"
$newZh  = "$newGuid.160268b711310859adf292a62b10ee58a83ab059.zh-cn.xlf"

$oldDe  = "$oldGuid.fddd2aa3a2d4333e1266330065820dd4a060f147.de-de.xlf"
$newDe  = "$newGuid.160268b711310859adf292a62b10ee58a83ab059.de-de.xlf"

$oldZhDate = "2016-03-09 12:53:23"
$newZhDate = "2016-03-09 12:53:59"

$oldDeDate = "2016-03-09 12:53:31"
$newDeDate = "2016-03-09 12:54:07"

# NOTE: only the *display* label of each hyperlink moves to the new GUID;
# the link target itself keeps pointing at the original (old-GUID) path, the
# report generator does not recompute the stored commit / target URLs.

# --- Sheet 1: "Overview" -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$mdTarget1 = "https://github.com/OpenLocalizationTest/oltest/blob/9b007c72fe3126f3d2fa830b6dbabe7c90a006e8/e2e/$oldMd"
$cfgTarget1 = "https://github.com/OpenLocalizationTest/oltest/blob/9b007c72fe3126f3d2fa830b6dbabe7c90a006e8/.localization-config"

$ws1.Range("A2").Value = $newMd

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdTarget1, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $cfgTarget1, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# --- Sheet 2: "zh-cn" -----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$mdTarget2  = "https://github.com/OpenLocalizationTest/oltest/blob/9b007c72fe3126f3d2fa830b6dbabe7c90a006e8/e2e/$oldMd"
$zhTarget2  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/335d8795e80843a73df7212602ec9de8ece68b90/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZh"
$cfgTarget2 = "https://github.com/OpenLocalizationTest/oltest/blob/9b007c72fe3126f3d2fa830b6dbabe7c90a006e8/.localization-config"

$ws2.Range("A2").Value = $newMd
$ws2.Range("C2").Value = $newZh
$ws2.Range("D2").Value = $newZhDate

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdTarget2, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhTarget2, [Type]::Missing, [Type]::Missing, $newZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $cfgTarget2, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# --- Sheet 3: "de-de" -----------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$mdTarget3  = "https://github.com/OpenLocalizationTest/oltest/blob/9b007c72fe3126f3d2fa830b6dbabe7c90a006e8/e2e/$oldMd"
$deTarget3  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfdff84bf64b47aebb1a0e15547bb48c386af2c1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDe"
$cfgTarget3 = "https://github.com/OpenLocalizationTest/oltest/blob/9b007c72fe3126f3d2fa830b6dbabe7c90a006e8/.localization-config"

$ws3.Range("A2").Value = $newMd
$ws3.Range("C2").Value = $newDe
$ws3.Range("D2").Value = $newDeDate

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdTarget3, [Type]::Missing, [Type]::Missing, $newMd) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deTarget3, [Type]::Missing, [Type]::Missing, $newDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $cfgTarget3, [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
